$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (column B) entirely, shifting
# columns C:H left to B:G.
$ws.Columns("B:B").Delete()

# Append ".global" to the (now shifted) header labels in B1:G1,
# leaving the "Country" header in A1 untouched.
foreach ($col in @("B", "C", "D", "E", "F", "G")) {
    $cell = $ws.Range("$col`1")
    $cell.Value = $cell.Value() + ".global"
}
